$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new row at row 9 (pushes old rows 9..67 down to 10..68).
#    This realigns every formula reference automatically (D column
#    running-balance formulas, the SUM(I3:I44) -> SUM(I3:I45) range).
# ------------------------------------------------------------------
$ws.Rows("9:9").Insert()

# ------------------------------------------------------------------
# 2) The freshly inserted row 9 comes back with a blank/derived style;
#    repair H9:I9 formatting by pulling it from row 10 (which holds
#    the style that used to belong to the old row 9), then fill in
#    the new "interest collected" line-item.
# ------------------------------------------------------------------
$ws.Range("H10:I10").Copy()
$ws.Range("H9:I9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("H9").Value = "21/02/2024"
$ws.Range("I9").Value = 10

# K9 gets the "chua lam giay" note flag - copy the format from the
# existing note cell (now at K19 after the insert) and reuse its
# exact text so the shared string is reused, not duplicated.
$ws.Range("K19").Copy()
$ws.Range("K9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("K9").Value = $ws.Range("K19").Value2

# ------------------------------------------------------------------
# 3) Append the two new ledger rows at the bottom of the A:D table
#    (new last rows 69 and 70). New cells automatically inherit the
#    formatting of the row above them.
# ------------------------------------------------------------------
$ws.Range("A69").Value = "21/02/2024"
$ws.Range("B69").Value = "Duy lấy tiền lời 5tr"
$ws.Range("C69").Value = 5000
$ws.Range("D69").Formula = "=D68+C69"

$ws.Range("B70").Value = $ws.Range("B11").Value2
$ws.Range("C70").Value = -10000
$ws.Range("D70").Formula = "=D69+C70"

# ------------------------------------------------------------------
# 4) Restore the selected cell shown in the saved view.
# ------------------------------------------------------------------
$ws.Range("K27").Select()
